$d = $word.ActiveDocument

# 1) "EN RT CERRAMOS CON ... VS LW ..." line: re-cased to "En RT cerramos con ... vs LW ..."
$d.Content.Find.Execute(
    "EN RT CERRAMOS CON {{rt_today}} VS LW {{rt_lw}}",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "En RT cerramos con {{rt_today}} vs LW {{rt_lw}}", 2) | Out-Null

# 2) "LA SUBCATEGORIA WO ES ..." line: re-worded to "LA razón Wo es ..."
$d.Content.Find.Execute(
    "LA SUBCATEGORIA WO ES {{catwosubrt}} con {{valorwosubrt}} vs LW {{valorwosubrtlw}}",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "LA razón Wo es {{catwosubrt}} con {{valorwosubrt}} vs LW {{valorwosubrtlw}}", 2) | Out-Null

# 3) "EN PARTNER CERRAMOS CON ... VS LW ..." line: re-cased to "En Partner cerramos con ... vs LW ..."
$d.Content.Find.Execute(
    "EN PARTNER CERRAMOS CON {{partner_today}} VS LW {{partner_lw}}",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "En Partner cerramos con {{partner_today}} vs LW {{partner_lw}}", 2) | Out-Null

# 4) Finish the partner sub-category line (the only one of the two identical
#    "La Razón Wo es" paragraphs that sits right after the Partner paragraph):
#    append the missing placeholders, mirroring the RT sub-category paragraph.
#    Scope the search to that single paragraph so the other, still-unfinished
#    "La Razón Wo es" paragraph (before "EN USER ...") is left untouched.
$partnerHeading = $d.Paragraphs.Item(5)
$subPartner = $partnerHeading.Next()
$subRange = $subPartner.Range
$subRange.MoveEnd(1, -1) | Out-Null
$subRange.Find.Execute(
    "La Razón Wo es",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "La Razón Wo es {{catwosubpa}} con {{valorwosubpa}} vs LW {{valorwosubpalw}}", 2) | Out-Null
